$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<...>_old" -> "<...>_FV2410", "<...>_new" -> "<...>_FV2504"
#    (the "diff" header in column K stays untouched)
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into a real Excel Table ("Table1") spanning the
#    full used range (A1:U85), headers included.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U85")
$tbl = $ws.ListObjects().Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
# The source workbook's header row already carries its own bold/fill
# formatting, so the table itself should not layer a named quick-style on
# top of it.
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, top-left of the frozen pane
#    is A2).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
